# COREINTL_holdings.xlsx update:
#  - bump the "Model holdings provided as of ..." date in the confidential
#    disclaimer from 2021-05-28 to 2021-06-09
#  - refresh the Weight / Percent Change figures for the two holdings rows
#    (EFA / EEM) and the Total row's Percent Change

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; temporarily unprotect so the cells can be
# updated, then restore protection afterwards.
$ws.Unprotect()

$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-09 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.8471053012596438
$ws.Range("E2").Value = -0.002685874740568917

$ws.Range("D3").Value = 0.1528946987403562
$ws.Range("E3").Value = -0.002886002886002825

$ws.Range("E4").Value = -0.00271647327307456

# The newline we just wrote into A7 makes the engine auto-size row 7;
# AutoFit it back down so the row keeps its original (implicit/default)
# height instead of persisting an explicit customHeight.
$ws.Rows(7).EntireRow.AutoFit()

$ws.Protect()
